# Updates the cryptos price list (Price + Volume(1h) columns) to the latest
# snapshot, and reorders two pairs of rows whose rank changed (Dai/PancakeSwap
# and NEARProtocol/FraxShare).

function Set-TextCell($ws, $addr, $val) {
    # Forces the cell to store $val as literal text, even when it looks like
    # a number (e.g. "42.960.86" or "1.00"), matching the inline-string cells
    # already used throughout this sheet. The original style is restored so
    # no unintended formatting changes are introduced.
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Per-row cell updates (Price / Volume(1h) columns) ---
Set-TextCell $ws "D2" "42.960.86"
$ws.Range("E2").Value = "  +0.26%  "

Set-TextCell $ws "D3" "2.364.11"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  +0.05%  "

Set-TextCell $ws "D5" "302.50"
$ws.Range("E5").Value = "  +0.38%  "

Set-TextCell $ws "D6" "95.72"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("E7").Value = "  +0.03%  "

Set-TextCell $ws "D8" "0.502"
$ws.Range("E8").Value = "  -0.38%  "

Set-TextCell $ws "D9" "0.490"
$ws.Range("E9").Value = "  -0.49%  "

Set-TextCell $ws "D10" "34.11"
$ws.Range("E10").Value = "  -0.24%  "

Set-TextCell $ws "D11" "0.125"
$ws.Range("E11").Value = "  +3.70%  "

$ws.Range("E12").Value = "  +0.14%  "

Set-TextCell $ws "D13" "18.39"
$ws.Range("E13").Value = "  -3.06%  "

$ws.Range("E14").Value = "  -0.14%  "

Set-TextCell $ws "D15" "2.730.64"
$ws.Range("E15").Value = "  +1.97%  "

Set-TextCell $ws "D16" "2.373.15"
$ws.Range("E16").Value = "  +1.05%  "

$ws.Range("E17").Value = "  +0.41%  "

Set-TextCell $ws "D18" "42.920.81"
$ws.Range("E18").Value = "  +0.39%  "

Set-TextCell $ws "D19" "11.93"
$ws.Range("E19").Value = "  -2.10%  "

$ws.Range("E20").Value = "  +1.63%  "

Set-TextCell $ws "D21" "0.0₃0886"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("E22").Value = "  +0.38%  "

Set-TextCell $ws "D23" "235.23"
$ws.Range("E23").Value = "  -0.15%  "

$ws.Range("E24").Value = "  -5.14%  "

Set-TextCell $ws "D27" "24.52"
$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("E28").Value = "  +0.81%  "

$ws.Range("E29").Value = "  +1.93%  "

Set-TextCell $ws "D30" "31.99"
$ws.Range("E30").Value = "  -0.70%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("E32").Value = "  +0.16%  "

Set-TextCell $ws "D33" "17.50"
$ws.Range("E33").Value = "  -1.06%  "

Set-TextCell $ws "D34" "131.41"
$ws.Range("E34").Value = "  -10.98%  "

Set-TextCell $ws "D35" "0.0715"
$ws.Range("E35").Value = "  +2.00%  "

Set-TextCell $ws "D36" "1.85"
$ws.Range("E36").Value = "  +2.65%  "

$ws.Range("E37").Value = "  +3.66%  "

Set-TextCell $ws "D38" "4.34"
$ws.Range("E38").Value = "  -2.01%  "

Set-TextCell $ws "D39" "2.83"
$ws.Range("E39").Value = "  +3.13%  "

$ws.Range("E40").Value = "  -2.07%  "

Set-TextCell $ws "D41" "0.108"
$ws.Range("E41").Value = "  -0.52%  "

Set-TextCell $ws "D42" "21.14"
$ws.Range("E42").Value = "  -3.66%  "

Set-TextCell $ws "D43" "1.932.25"
$ws.Range("E43").Value = "  +0.53%  "

$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("E45").Value = "  +2.77%  "

Set-TextCell $ws "D48" "2.595.86"
$ws.Range("E48").Value = "  +1.92%  "

$ws.Range("E49").Value = "  +2.32%  "

$ws.Range("E50").Value = "  +1.63%  "

Set-TextCell $ws "D51" "71.47"

# --- Row swap: Dai (row 25) and PancakeSwap (row 26) traded ranking ---
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextCell $ws "D25" "2.43"
$ws.Range("E25").Value = "  +0.78%  "

$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextCell $ws "D26" "1.00"
$ws.Range("E26").Value = "  -0.10%  "

# --- Row swap: NEARProtocol (row 46) and FraxShare (row 47) traded ranking ---
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextCell $ws "D46" "9.18"
$ws.Range("E46").Value = "  -8.95%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextCell $ws "D47" "2.71"
$ws.Range("E47").Value = "  -1.17%  "
